$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.340.57"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "2.062.17"
$ws.Range("E3").Value = "  +4.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'237.15"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("D7").Value = "'58.13"
$ws.Range("E7").Value = "  +7.94%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +3.99%  "
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'0.0760"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "'0.102"
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("D13").Value = "2.366.61"
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").Value = "'14.53"
$ws.Range("E14").Value = "  +5.30%  "
$ws.Range("D15").Value = "'21.05"
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("D16").Value = "'0.778"
$ws.Range("E16").Value = "  +4.28%  "
$ws.Range("D17").Value = "'5.18"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").Value = "2.063.50"
$ws.Range("E18").Value = "  +4.64%  "
$ws.Range("D19").Value = "37.572.76"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").Value = "'6.14"
$ws.Range("E20").Value = "  +17.87%  "
$ws.Range("D21").Value = "'69.20"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").Value = "0.0₃0814"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'225.77"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "'163.89"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "'1.48"
$ws.Range("E28").Value = "  +12.23%  "
$ws.Range("D29").Value = "'8.86"
$ws.Range("E29").Value = "  +4.33%  "
$ws.Range("D30").Value = "'19.12"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "'4.49"
$ws.Range("E33").Value = "  +3.79%  "
$ws.Range("D34").Value = "'0.0625"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").Value = "'2.58"
$ws.Range("E35").Value = "  +12.38%  "
$ws.Range("D36").Value = "'4.49"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.37"
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "'5.87"
$ws.Range("E40").Value = "  +11.40%  "
$ws.Range("D41").Value = "'0.0991"
$ws.Range("E41").Value = "  +12.51%  "
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("E43").Value = "  +26.12%  "
$ws.Range("D44").Value = "'97.36"
$ws.Range("E44").Value = "  +10.98%  "
$ws.Range("D45").Value = "1.479.15"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("E46").Value = "  +7.61%  "
$ws.Range("E47").Value = "  +5.33%  "
$ws.Range("D48").Value = "'15.91"
$ws.Range("E48").Value = "  +8.33%  "
$ws.Range("D49").Value = "'1.03"
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("D50").Value = "'7.20"
$ws.Range("E50").Value = "  +7.20%  "
$ws.Range("E51").Value = "  +2.67%  "
